$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-33 from 2023-10-22 (serial 45221)
# to 2023-10-25 (serial 45224), keeping the value a clean date (no time component).
$newDateSerial = 45224

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = $newDateSerial
    }
}
